# Apply the "Summer Institute Edition '14" -> "Spring 2015" title update,
# and drop the unused "Footnote Text" paragraph style from styles.xml.

$d = $word.ActiveDocument

# --- 1. Update the course/term line on the title page -----------------
# Use a direct Range.Text assignment (rather than Find/Replace) so the
# existing run formatting (bold) and the xml:space="preserve" attribute
# on the <w:t> element are left intact, matching how Word itself edits
# text in place.
$target = "Educational Technology 0858-501, Summer Institute Edition " + [char]8217 + "14"
$replacement = "Educational Technology 0858-501, Spring 2015"

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $target) {
        $p.Range.Text = $replacement
        break
    }
}

# --- 2. Remove the unused "Footnote Text" style ------------------------
# It is not referenced anywhere in the document body or footnotes, so it
# can simply be deleted from the style sheet.
for ($i = $d.Styles.Count; $i -ge 1; $i--) {
    $s = $d.Styles.Item($i)
    if (($s.NameLocal -eq "Footnote Text") -or ($s.NameLocal -eq "FootnoteText")) {
        $s.Delete()
    }
}
